# Refresh the cryptocurrency market snapshot (prices in column D, 1h
# volume change in column E) with the latest pull from coinranking.com.
# Cardano/Dogecoin also swap rows 8/9 back to their canonical rank order.
#
# Numeric-looking price strings are entered with a leading apostrophe so
# Excel keeps them as text (matching the sheet's existing text-formatted
# price column) instead of auto-converting to a General number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.381.47"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.875.99"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D5").Value = "'0.7118"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'242.12"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.07798"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3115"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").Value = "'25.21"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").Value = "'0.08433"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").Value = "1.865.04"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "'0.7143"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "'91.10"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "29.381.20"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'6.061"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("E18").Value = "  +5.02%  "
$ws.Range("D19").Value = "'240.69"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "2.120.68"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'7.783"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "'163.26"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "'9.068"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("D28").Value = "'18.56"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "'4.427"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'4.329"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("D33").Value = "'0.05293"
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").Value = "'1.942"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "'0.7458"
$ws.Range("E36").Value = "  -9.88%  "
$ws.Range("D37").Value = "'2.695"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").Value = "'0.01871"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "1.230.81"
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("D40").Value = "'2.730"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("E42").Value = "  +8.13%  "
$ws.Range("D43").Value = "'0.8912"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "'72.73"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'0.9999"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "2.017.46"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "'1.814"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").Value = "'0.5213"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").Value = "'9.419"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'0.4324"
$ws.Range("E51").Value = "  +1.23%  "
